$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7452577948570251
$ws.Range("B1").Value = 1.4085533618927
$ws.Range("D1").Value = 1.923132061958313
$ws.Range("E1").Value = 1.213793635368347
